$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.448.13'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.371.25'
$ws.Range('E3').Value = '  +2.53%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.77'
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.91'
$ws.Range('E6').Value = '  +2.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.509'
$ws.Range('E7').Value = '  -4.93%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.519'
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.66'
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.25'
$ws.Range('E11').Value = '  +1.96%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0805'
$ws.Range('E12').Value = '  -1.46%  '
$ws.Range('E13').Value = '  -0.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.93'
$ws.Range('E14').Value = '  -3.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.51'
$ws.Range('E15').Value = '  +3.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.372.64'
$ws.Range('E16').Value = '  +2.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.809'
$ws.Range('E17').Value = '  -0.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.435.90'
$ws.Range('E18').Value = '  +0.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.32'
$ws.Range('E19').Value = '  +3.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.88'
$ws.Range('E20').Value = '  -5.54%  '
$ws.Range('E21').Value = '  -0.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.13'
$ws.Range('E22').Value = '  -0.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.10'
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('E24').Value = '  +0.60%  '
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.81'
$ws.Range('E27').Value = '  +3.94%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.32'
$ws.Range('E28').Value = '  +10.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.44'
$ws.Range('E29').Value = '  -2.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.46'
$ws.Range('E30').Value = '  -2.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '161.54'
$ws.Range('E31').Value = '  -2.98%  '
$ws.Range('E32').Value = '  -2.40%  '
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.14'
$ws.Range('E34').Value = '  +0.76%  '
$ws.Range('E35').Value = '  +5.06%  '
$ws.Range('E36').Value = '  -2.66%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.66'
$ws.Range('E37').Value = '  +8.17%  '
$ws.Range('E38').Value = '  -1.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.91'
$ws.Range('E39').Value = '  +3.58%  '
$ws.Range('E40').Value = '  -2.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.113'
$ws.Range('E41').Value = '  -2.31%  '
$ws.Range('E42').Value = '  +12.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.031.72'
$ws.Range('E43').Value = '  +2.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.57'
$ws.Range('E44').Value = '  +0.53%  '
$ws.Range('E45').Value = '  -0.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.52'
$ws.Range('E46').Value = '  +6.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.09'
$ws.Range('E47').Value = '  +2.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '57.58'
$ws.Range('E48').Value = '  +3.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.92'
$ws.Range('E49').Value = '  -1.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.607.88'
$ws.Range('E50').Value = '  +2.64%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.79'
$ws.Range('E51').Value = '  +7.23%  '
